$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.735.57'
$ws.Range('E2').Value = '  +2.36%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.311.76'
$ws.Range('E3').Value = '  +1.13%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '319.73'
$ws.Range('E5').Value = '  +1.61%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '104.83'
$ws.Range('E6').Value = '  +2.25%  '
$ws.Range('E7').Value = '  +1.01%  '
$ws.Range('E8').Value = '  +0.19%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.611'
$ws.Range('E9').Value = '  +1.84%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.32'
$ws.Range('E10').Value = '  +3.91%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0911'
$ws.Range('E11').Value = '  +0.88%  '
$ws.Range('E12').Value = '  +5.09%  '
$ws.Range('E13').Value = '  +1.12%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.981'
$ws.Range('E14').Value = '  +2.31%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.43'
$ws.Range('E15').Value = '  +1.21%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.656.81'
$ws.Range('E16').Value = '  +1.14%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.329.19'
$ws.Range('E17').Value = '  +1.71%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '42.846.95'
$ws.Range('E18').Value = '  +2.74%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.56'
$ws.Range('E19').Value = '  +2.00%  '
$ws.Range('E20').Value = '  +1.37%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.58'
$ws.Range('E21').Value = '  +34.61%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '73.97'
$ws.Range('E22').Value = '  +1.32%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.57'
$ws.Range('E23').Value = '  -1.43%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '271.90'
$ws.Range('E24').Value = '  -2.69%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.28'
$ws.Range('E25').Value = '  +1.64%  '
$ws.Range('E26').Value = '  -0.48%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.97'
$ws.Range('E27').Value = '  +2.66%  '
$ws.Range('E28').Value = '  -3.17%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '22.76'
$ws.Range('E29').Value = '  -0.29%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '38.18'
$ws.Range('E30').Value = '  +9.19%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.29'
$ws.Range('E31').Value = '  +8.66%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '166.08'
$ws.Range('E32').Value = '  +1.78%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0894'
$ws.Range('E33').Value = '  +2.92%  '
$ws.Range('E34').Value = '  +0.78%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.116'
$ws.Range('E35').Value = '  +0.72%  '
$ws.Range('E36').Value = '  -10.72%  '
$ws.Range('E37').Value = '  +2.35%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0356'
$ws.Range('E38').Value = '  +3.14%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.72'
$ws.Range('E39').Value = '  +3.00%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.75'
$ws.Range('E40').Value = '  -4.42%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.58'
$ws.Range('E41').Value = '  +9.16%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '100.87'
$ws.Range('E42').Value = '  +1.26%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '70.85'
$ws.Range('E43').Value = '  +2.27%  '
$ws.Range('E44').Value = '  +1.86%  '
$ws.Range('E45').Value = '  +0.23%  '
$ws.Range('B46').Value = 'ordi'
$ws.Range('C46').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '83.33'
$ws.Range('E46').Value = '  +10.19%  '
$ws.Range('B47').Value = 'Celestia'
$ws.Range('C47').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '12.40'
$ws.Range('E47').Value = '  +4.96%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '114.99'
$ws.Range('E48').Value = '  -0.52%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.33'
$ws.Range('E49').Value = '  +1.92%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.91'
$ws.Range('E50').Value = '  +0.05%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.601.66'
$ws.Range('E51').Value = '  +5.09%  '
